# Edit: append two trailing spaces to the first paragraph's sentence, then
# add a new run right after it (same paragraph) containing a dark-red
# (C00000) note about an alternate branch version.

$d = $word.ActiveDocument

# The first paragraph currently reads exactly:
#   "This is a Microsoft word document."
$p1 = $d.Paragraphs(1).Range

# Build a range that covers the paragraph's text but excludes the trailing
# paragraph mark, so InsertAfter appends inside the paragraph (not after it).
$r = $p1.Duplicate
$r.MoveEnd(1, -1) | Out-Null   # wdCharacter = 1; shrink end by one char (the pilcrow)

# 1) Add two trailing spaces after the existing sentence, still as part of
#    the original (unformatted) run.
$r.InsertAfter("  ")

# 2) Collapse the range to its new end and insert a second run with the
#    red "(This is a change - Version for branch alternate)" annotation.
$r.Collapse(0)      # wdCollapseEnd = 0
$r.InsertAfter("(This is a change – Version for branch alternate)")

# Word's Font.Color is a BGR-packed long (R + G*256 + B*65536). For
# RRGGBB = C00000 (R=0xC0, G=0x00, B=0x00) that packed value is simply
# 0xC0 = 192.
$r.Font.Color = 192
